# "new Madigan bike hours" - updated weekly ridership numbers.
# Columns: A=Weekday, B=Date, C=Riders, D=Average, E=Pilot Target.
# Only the Riders (C) and Average (D) columns change for this week.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# row -> (Riders, Average)
$updates = @(
    @{ Row = 2;  Riders = 277; Average = 233.47 },
    @{ Row = 3;  Riders = 219; Average = 212.9 },
    @{ Row = 4;  Riders = 215; Average = 211.9 },
    @{ Row = 5;  Riders = 204; Average = 238.1 },
    @{ Row = 6;  Riders = 305; Average = 240.71 },
    @{ Row = 7;  Riders = 85;  Average = 115.45 },
    @{ Row = 8;  Riders = 69;  Average = 95.25 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.Riders    # column C - Riders
    $ws.Cells.Item($u.Row, 4).Value = $u.Average   # column D - Average
}

# Best-effort refresh of the embedded "GO Transit Ridership" chart so it
# picks up the new Riders/Average series values.
try {
    $co = $ws.ChartObjects(1)
    $co.Chart.Refresh()
} catch {
}

$wb.Save()
